$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 21:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 390387
$ws.Range("C4").Value = 23383
$ws.Range("D4").Value = 21488
$ws.Range("E4").Value = 356425
$ws.Range("F4").Value = 9169
$ws.Range("G4").Value = 1603
$ws.Range("H4").Value = 12474

# Row 49 - Republica Dominicana
$ws.Range("D49").Value = 36
$ws.Range("E49").Value = 1822

# Row 77 - Camerun
$ws.Range("D77").Value = 43
$ws.Range("E77").Value = 606

# Row 84 - Principado de Andorra
$ws.Range("B84").Value = 545
$ws.Range("C84").Value = 20
$ws.Range("D84").Value = 39
$ws.Range("E84").Value = 484
$ws.Range("F84").Value = 17
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 22

# Row 119 - Martinica
$ws.Range("B119").Value = 152
$ws.Range("C119").Value = 1
$ws.Range("E119").Value = 98

# Row 130 - Madagascar
$ws.Range("B130").Value = 88
$ws.Range("C130").Value = 6
$ws.Range("E130").Value = 86

# Row 152 - Bahamas
$ws.Range("E152").Value = 22
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 6
